# Update "想去人数" (interested-count) values in column F across sheets,
# matching the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1664
$ws1.Range("F3").Value = 9082
$ws1.Range("F5").Value = 506
$ws1.Range("F6").Value = 702
$ws1.Range("F7").Value = 1323
$ws1.Range("F8").Value = 197
$ws1.Range("F11").Value = 5887
$ws1.Range("F13").Value = 385
$ws1.Range("F15").Value = 4393
$ws1.Range("F16").Value = 13
$ws1.Range("F17").Value = 163
$ws1.Range("F19").Value = 24
$ws1.Range("F21").Value = 19
$ws1.Range("F22").Value = 254
$ws1.Range("F23").Value = 13
$ws1.Range("F24").Value = 2729

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 22

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1664
$ws4.Range("F3").Value = 9082
$ws4.Range("F5").Value = 22
$ws4.Range("F6").Value = 506
$ws4.Range("F7").Value = 702
$ws4.Range("F8").Value = 1323
$ws4.Range("F9").Value = 197
$ws4.Range("F12").Value = 5887
$ws4.Range("F14").Value = 385
$ws4.Range("F16").Value = 4393
$ws4.Range("F17").Value = 13
$ws4.Range("F18").Value = 163
$ws4.Range("F20").Value = 24
$ws4.Range("F22").Value = 19
$ws4.Range("F23").Value = 254
$ws4.Range("F24").Value = 13
$ws4.Range("F25").Value = 2729
